$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number need a leading
# apostrophe so Excel stores them as text (matching the original inlineStr
# cells) instead of converting them to a numeric cell. ClearFormats() afterwards
# removes the quote-prefix cell style so no stray formatting is introduced.

$ws.Range("D2").Value = '61.633.10'
$ws.Range("E2").Value = '  -3.95%  '
$ws.Range("D3").Value = '2.975.68'
$ws.Range("E3").Value = '  -5.01%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''543.57'
$ws.Range("E5").Value = '  -4.66%  '
$ws.Range("D6").Value = '''152.61'
$ws.Range("E6").Value = '  -5.28%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '''0.575'
$ws.Range("E8").Value = '  +1.83%  '
$ws.Range("D9").Value = '2.984.89'
$ws.Range("E9").Value = '  -4.90%  '
$ws.Range("E10").Value = '  -1.60%  '
$ws.Range("D11").Value = '''6.15'
$ws.Range("E11").Value = '  -6.23%  '
$ws.Range("E12").Value = '  -1.79%  '
$ws.Range("D13").Value = '3.497.29'
$ws.Range("E13").Value = '  -4.82%  '
$ws.Range("E14").Value = '  -2.31%  '
$ws.Range("D15").Value = '61.681.75'
$ws.Range("E15").Value = '  -3.91%  '
$ws.Range("D16").Value = '''23.75'
$ws.Range("E16").Value = '  -4.31%  '
$ws.Range("D17").Value = '2.985.12'
$ws.Range("E17").Value = '  -4.81%  '
$ws.Range("D18").Value = '''0.0000147'
$ws.Range("E18").Value = '  -3.62%  '
$ws.Range("D19").Value = '''5.18'
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("D20").Value = '''382.37'
$ws.Range("E20").Value = '  -4.46%  '
$ws.Range("D21").Value = '''12.03'
$ws.Range("E21").Value = '  -3.56%  '
$ws.Range("D22").Value = '''6.69'
$ws.Range("E22").Value = '  -5.46%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").Value = '''65.82'
$ws.Range("E24").Value = '  -2.94%  '
$ws.Range("E25").Value = '  -1.67%  '
$ws.Range("D26").Value = '3.101.00'
$ws.Range("E26").Value = '  -5.33%  '
$ws.Range("D27").Value = '''0.190'
$ws.Range("E27").Value = '  -2.02%  '
$ws.Range("D28").Value = '0.0₃0944'
$ws.Range("E28").Value = '  -5.51%  '
$ws.Range("D29").Value = '''0.994'
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("D30").Value = '''8.30'
$ws.Range("E30").Value = '  -5.40%  '
$ws.Range("D31").Value = '''0.999'
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").Value = '''1.72'
$ws.Range("E32").Value = '  -4.14%  '
$ws.Range("D33").Value = '''20.50'
$ws.Range("E33").Value = '  -2.75%  '
$ws.Range("D34").Value = '''160.67'
$ws.Range("E34").Value = '  +1.44%  '
$ws.Range("D35").Value = '''4.71'
$ws.Range("E35").Value = '  -1.40%  '
$ws.Range("D36").Value = '''5.97'
$ws.Range("E36").Value = '  -4.00%  '
$ws.Range("E37").Value = '  -1.84%  '
$ws.Range("D38").Value = '''1.27'
$ws.Range("E38").Value = '  -4.16%  '
$ws.Range("D39").Value = '''1.56'
$ws.Range("E39").Value = '  -5.64%  '
$ws.Range("D40").Value = '''3.94'
$ws.Range("E40").Value = '  -2.79%  '
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").Value = '''37.40'
$ws.Range("E41").Value = '  -2.42%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.410.78'
$ws.Range("E42").Value = '  -9.31%  '
$ws.Range("D43").Value = '''22.20'
$ws.Range("E43").Value = '  -4.90%  '
$ws.Range("E44").Value = '  -2.33%  '
$ws.Range("D45").Value = '''0.0595'
$ws.Range("E45").Value = '  -2.08%  '
$ws.Range("D46").Value = '''0.0248'
$ws.Range("E46").Value = '  -2.00%  '
$ws.Range("D47").Value = '''5.09'
$ws.Range("E47").Value = '  -5.77%  '
$ws.Range("D48").Value = '''0.997'
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("D49").Value = '''270.81'
$ws.Range("E49").Value = '  -6.00%  '
$ws.Range("D50").Value = '''19.89'
$ws.Range("E50").Value = '  -5.03%  '
$ws.Range("D51").Value = '''0.0955'
$ws.Range("E51").Value = '  -1.56%  '

# Strip the quote-prefix style iron_native applied to the text-forced cells above
# so their cell style matches the rest of the (unstyled) data cells.
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
